$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Executives")
$ws.Activate()

# Delete the entire AA column (a stray duplicate column), shifting everything
# after it one column to the left.
$ws.Columns("AA").Delete()

# Fix up the header row labels in the now-shifted block (U3:Z3) so the
# Shares/Price/Perf RS Value/Shares/Price/Perf Cash sequence reads correctly.
$ws.Range("U3").Value = "Shares"
$ws.Range("V3").Value = "Price"
$ws.Range("W3").Value = "Perf RS Value"
$ws.Range("X3").Value = "Shares"
$ws.Range("Y3").Value = "Price"
$ws.Range("Z3").Value = "Perf Cash"
